$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.471.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.916.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.46%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4749"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.92%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4090"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08054"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.24%  "

$ws.Range("E11").Value = "  +1.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.861.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.925"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.153"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.52%  "

$ws.Range("E17").Value = "  -0.72%  "

$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.73"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.71%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.487.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.539"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.207"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.102.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.890"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.132"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.31%  "

$ws.Range("E32").Value = "  +7.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09541"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.429"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.567"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.400"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.91%  "

$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.313"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.174"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5884"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.543"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1844"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.68%  "

$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.08002"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.65%  "

$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.19"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.34%  "

$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5545"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.28%  "

$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.24%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.99%  "
